$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns for zh-cn (E) and de-de (F) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Column widths widen to fit the longer status text (auto-fit side effect)
$wsOverview.Columns.Item(5).ColumnWidth = 29.09
$wsOverview.Columns.Item(6).ColumnWidth = 29.09

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# Latest Handback DateTime refreshed for the new handback run
$wsZhCn.Range("K2").Value = "2016-08-31 12:09:04"
$wsZhCn.Range("K3").Value = "2016-08-31 12:09:04"

# Error Detail cleared now that the handback is in sync (keep as an empty text cell)
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = $wsZhCn.Range("B2").Style
$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Range("P3").Style = $wsZhCn.Range("B3").Style

$wsZhCn.Columns.Item(3).ColumnWidth = 29.09
$wsZhCn.Columns.Item(16).ColumnWidth = 12.75

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Range("K2").Value = "2016-08-31 12:09:24"
$wsDeDe.Range("K3").Value = "2016-08-31 12:09:24"

$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = $wsDeDe.Range("B2").Style
$wsDeDe.Range("P3").Value = "'"
$wsDeDe.Range("P3").Style = $wsDeDe.Range("B3").Style

$wsDeDe.Columns.Item(3).ColumnWidth = 29.09
$wsDeDe.Columns.Item(16).ColumnWidth = 12.75
